$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fees")

# B23 date moved from 12/7/2020 (44172) to 12/8/2020 (44173)
$ws.Range("B23").Value = 44173

# Row 24: new invoice line (date, task, hours, cost formula)
$ws.Range("D23").Copy()
$ws.Range("D24:D25").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C24:C25").PasteSpecial(-4122)

$ws.Range("B24").Value = 44174
$ws.Range("C24").Value = "Clean up 2013 8/10th MTF data for FA and figure out NA issue"
$ws.Range("D24").Value = 2
$ws.Range("E24").Formula = "=D24*40"

# Row 25: new invoice line (date, task, hours, cost formula)
$ws.Range("B25").Value = 44175
$ws.Range("C25").Value = "Clean ALL years MTF data for FA "
$ws.Range("D25").Value = 1
$ws.Range("E25").Formula = "=D25*40"

$excel.Calculate()

# Move the saved selection/active cell to G22
$ws.Activate()
$ws.Range("G22").Select()
